# Update schedule to emphasize R in initial weeks.
# Rows 3-12 in column D (the "topic" text for cm001-cm010) are rewritten so
# the curriculum introduces R earlier; rows 13-22 (Web scraping ... HPC)
# keep their existing text. The active selection also moves to D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = "Introduction to computational social science, R, git, and workflow"
$ws.Range("D4").Value  = "Data visualization"
$ws.Range("D5").Value  = "Exploratory data analysis"
$ws.Range("D6").Value  = "Data wrangling"
$ws.Range("D7").Value  = "Programming in R"
$ws.Range("D8").Value  = "Programming in R"
$ws.Range("D9").Value  = "Models in R"
$ws.Range("D10").Value = "Introduction to Python"
$ws.Range("D11").Value = "Functional programming in Python"
$ws.Range("D12").Value = "Data pipelines, communication, and reproducibility"

# Move the sheet's active selection from A2 to D13.
$ws.Range("D13").Select() | Out-Null
